$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Department")

$ws.Range("A25").Value = "158 W 124TH ST IAM DENTAL"
$ws.Range("B25").Value = 8274005
$ws.Range("C25").Value = "MSM"

$ws.Range("A26").Value = "440 W 114TH IAM DENTAL"
$ws.Range("B26").Value = 8808068
$ws.Range("C26").Value = "MSM"

$ws.Range("A27").Value = "X_440 W 114TH IAM DENTAL_DEACTIVATED"
$ws.Range("B27").Value = 43006057
$ws.Range("C27").Value = "MSM"

$ws.Range("B23").Copy()
$ws.Range("A25:B27").PasteSpecial(-4122)

$ws.Range("G29").Select() | Out-Null

